$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The merged header I2:K2 now describes a single combined parcel ("3227-0")
# instead of the old "734-18" table.
$ws.Range("I2").Value = "Koordinatat e parceles 3227-0"

# Rows 5-7: the coordinate rows shift up by one (the old row 4 entry for
# parcel 168621 stays, the rest of the first block is replaced with values
# that used to belong to the second ("733-20") table).
$ws.Range("I5").Value = 168622
$ws.Range("J5").Value = 7511867.739
$ws.Range("K5").Value = 4692128.479

$ws.Range("I6").Value = 181505
$ws.Range("J6").Value = 7511887.445
$ws.Range("K6").Value = 4692119.66

$ws.Range("I7").Value = 181506
$ws.Range("J7").Value = 7511888.297
$ws.Range("K7").Value = 4692127.302

# Row 8 had no I:K values before; copy number formatting from row 4 so the
# new cells pick up the existing styles instead of creating new ones.
$ws.Range("I4:K4").Copy()
$ws.Range("I8:K8").PasteSpecial(-4122)
$ws.Range("I8").Value = 181509
$ws.Range("J8").Value = 7511886.701
$ws.Range("K8").Value = 4692113

# Row 9 used to hold the merged "Koordinatat e parceles 733-20" sub-header;
# unmerge it and turn it into a normal data row.
$ws.Range("I9:K9").UnMerge()
$ws.Range("I4:K4").Copy()
$ws.Range("I9:K9").PasteSpecial(-4122)
$ws.Range("I9").Value = 182181
$ws.Range("J9").Value = 7511866.27
$ws.Range("K9").Value = 4692115.73

$excel.CutCopyMode = 0

# Rows 10-16 used to hold the rest of the now-removed second table; clear
# them entirely so the <c> elements disappear from the sheet.
$ws.Range("I10:K16").Clear()

# Restore the on-screen selection to match the saved state.
$null = $ws.Range("V16").Select()
